# Update countries & provincias Spain
# Refresh the "Pais" data feed: new totals for a handful of countries plus
# re-ranking of three countries whose case totals overtook their neighbours,
# and bump the "datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 06:36"

# --- Straight data refresh (no re-ranking) -----------------------------

# Brasil (row 5)
$ws.Cells.Item(5, 4).Value = 916147
$ws.Cells.Item(5, 5).Value = 476509

# Pakistan (row 15)
$ws.Cells.Item(15, 2).Value = 217809
$ws.Cells.Item(15, 3).Value = 4339
$ws.Cells.Item(15, 4).Value = 104694
$ws.Cells.Item(15, 5).Value = 108642
$ws.Cells.Item(15, 7).Value = 78
$ws.Cells.Item(15, 8).Value = 4473

# Mongolia (row 169)
$ws.Cells.Item(169, 4).Value = 177
$ws.Cells.Item(169, 5).Value = 43

# Butan (row 186)
$ws.Cells.Item(186, 4).Value = 50
$ws.Cells.Item(186, 5).Value = 27

# --- Kazajistan overtakes Portugal and Oman (rows 39-41) ---------------
# Row 39 becomes Kazajistan with its updated totals
$ws.Cells.Item(39, 1).Value = "Kazajistan"
$ws.Cells.Item(39, 2).Value = 42574
$ws.Cells.Item(39, 3).Value = 1509
$ws.Cells.Item(39, 4).Value = 25319
$ws.Cells.Item(39, 5).Value = 17067
$ws.Cells.Item(39, 8).Value = 188

# Row 40 becomes Portugal (shifted down one position, values unchanged)
$ws.Cells.Item(40, 1).Value = "Portugal"
$ws.Cells.Item(40, 2).Value = 42454
$ws.Cells.Item(40, 4).Value = 27798
$ws.Cells.Item(40, 5).Value = 13077
$ws.Cells.Item(40, 8).Value = 1579

# Row 41 becomes Oman (shifted down one position, values unchanged)
$ws.Cells.Item(41, 1).Value = "Oman"
$ws.Cells.Item(41, 2).Value = 41194
$ws.Cells.Item(41, 4).Value = 24162
$ws.Cells.Item(41, 5).Value = 16847
$ws.Cells.Item(41, 8).Value = 185

# --- Kirguistan overtakes Venezuela, Haiti, Tayikistan, Etiopia (rows 82-86)

# Row 82 becomes Kirguistan with its updated totals
$ws.Cells.Item(82, 1).Value = "Kirguistan"
$ws.Cells.Item(82, 2).Value = 6261
$ws.Cells.Item(82, 3).Value = 526
$ws.Cells.Item(82, 4).Value = 2530
$ws.Cells.Item(82, 5).Value = 3665
$ws.Cells.Item(82, 7).Value = 4
$ws.Cells.Item(82, 8).Value = 66

# Row 83 becomes Venezuela (shifted down one position, values unchanged)
$ws.Cells.Item(83, 1).Value = "Venezuela"
$ws.Cells.Item(83, 2).Value = 6062
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 1649
$ws.Cells.Item(83, 5).Value = 4359
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 54

# Row 84 becomes Haiti (shifted down one position, with its own updated totals)
$ws.Cells.Item(84, 1).Value = "Haiti"
$ws.Cells.Item(84, 2).Value = 6040
$ws.Cells.Item(84, 3).Value = 65
$ws.Cells.Item(84, 4).Value = 1032
$ws.Cells.Item(84, 5).Value = 4901
$ws.Cells.Item(84, 7).Value = 2
$ws.Cells.Item(84, 8).Value = 107

# Row 85 becomes Tayikistan (shifted down one position, values unchanged)
$ws.Cells.Item(85, 1).Value = "Tayikistan"
$ws.Cells.Item(85, 2).Value = 6005
$ws.Cells.Item(85, 4).Value = 4627
$ws.Cells.Item(85, 5).Value = 1326
$ws.Cells.Item(85, 8).Value = 52

# Row 86 becomes Etiopia (shifted down one position, values unchanged)
$ws.Cells.Item(86, 1).Value = "Etiopia"
$ws.Cells.Item(86, 2).Value = 5846
$ws.Cells.Item(86, 4).Value = 2430
$ws.Cells.Item(86, 5).Value = 3313
$ws.Cells.Item(86, 8).Value = 103

# --- Santa Lucia overtakes Laos (identical totals, rows 203-204) -------
$ws.Cells.Item(203, 1).Value = "Santa Lucia"
$ws.Cells.Item(204, 1).Value = "Laos"
